$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: logistic_embeddings
$ws.Range("C5").Value = 0.183
$ws.Range("D5").Value = 0.343
$ws.Range("E5").Value = 0.4
$ws.Range("F5").Value = 0.458
$ws.Range("G5").Value = 0.492
$ws.Range("H5").Value = 0.513

# Row 7: classical-best-embeddings -> classical-best-embed
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.183
$ws.Range("D7").Value = 0.343
$ws.Range("E7").Value = 0.4

# Row 8: BERT-base
$ws.Range("C8").Value = 0.154
$ws.Range("D8").Value = 0.427
$ws.Range("E8").Value = 0.478
$ws.Range("F8").Value = 0.545
$ws.Range("G8").Value = 0.588
$ws.Range("H8").Value = 0.604

# Row 9: BERT-base-nli
$ws.Range("B9").Value = 0.271
$ws.Range("C9").Value = 0.388
$ws.Range("D9").Value = 0.484
$ws.Range("E9").Value = 0.51
$ws.Range("G9").Value = 0.5600000000000001
$ws.Range("H9").Value = 0.576
